# PNAD 2009 - correção nos dados e inicio da analise
# The source data had two "section header" rows (row 5 "situação do
# domicílio" and, two rows below it, "grandes regiões e unidades da
# federação") that carried no numeric values of their own - the real
# numbers for "urbana"/"rural"/"norte"/etc. had accidentally been written
# one row too low. This removes those two empty label rows (which shifts
# every following row up so the labels now line up with their correct
# data) and relabels the B2 header from the stray pandas artifact
# "unnamed: 1_level_1" to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-imported pandas column header.
$ws.Range("B2").Value = "total"

# Remove the two empty "section header" rows so the data shifts up into
# the correct label rows. Row 5 ("situação do domicílio") goes first;
# after it is removed, the old row 8 ("grandes regiões e unidades da
# federação") becomes row 7.
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(7).Delete() | Out-Null
